$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the filename shared string (affects H2, H3, and H4 via shared string table)
$ws.Range("H2").Value = "documents_true_table.json"
$ws.Range("H3").Value = "documents_true_table.json"
$ws.Range("H4").Value = "documents_true_table.json"

# Update quantity_of_rules (column E) and quantity_of_terms (column G)
$ws.Range("E2").Value = 2
$ws.Range("G2").Value = 50

$ws.Range("E3").Value = 4
$ws.Range("G3").Value = 36
